# Add season-record columns (Wins / Losses / Ties) to the roster table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers in row 1, copying the existing bold/bordered header style from AC1
# so the new cells pick up the same style index instead of creating a new one.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every player row (2-44).
$wins = 94
$losses = 68
$ties = 0

for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}

Write-Output "Season record columns (AD:AF) added for rows 1-44"
